# Added const for yield
# Fill in the "YIELD" column (R) with the constant value 345 for every
# data row that currently has an empty inline-string placeholder cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$yieldRange = $ws.Range("R5:R11,R16:R24,R29:R39,R44:R54,R59:R63,R68:R72,R77:R81,R86:R90,R95:R99,R104:R108,R113:R117")

foreach ($area in $yieldRange.Areas) {
    $area.NumberFormat = "@"
    $area.Value = "345"
}
